# Automatic update of files.
# Update the "Förändrad" (Changed) date column (C) for rows 2-6 from
# serial date 45174 (2023-09-05) to 45175 (2023-09-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = 45175
}
